# Auto-generated script applying numeric updates to Typhon_Profits market-data sheets
# (ALC, BSM, CRP, CUL, GSM, LTW) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 459.89474
$ws.Range("I33").Value = 402.25
$ws.Range("J33").Value = 767.3333
$ws.Range("K33").Value = 402.25
$ws.Range("L33").Value = 767.3333
$ws.Range("M33").Value = -173.25
$ws.Range("N33").Value = -1225.3333

$ws.Range("H100").Value = 2336.5789
$ws.Range("I100").Value = 1599.5454
$ws.Range("J100").Value = 3350
$ws.Range("K100").Value = 1599.5454
$ws.Range("L100").Value = 3350
$ws.Range("M100").Value = -1058.5454
$ws.Range("N100").Value = -4432

$ws.Range("H103").Value = 833546.7
$ws.Range("I103").Value = 1250070
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 3750210
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = -3749624
$ws.Range("N103").Value = -2672

$ws.Range("H133").Value = 51747.5
$ws.Range("J133").Value = 51747.5
$ws.Range("L133").Value = 51747.5
$ws.Range("N133").Value = -61867.5

$ws.Range("H135").Value = 27781824
$ws.Range("J135").Value = 83343000
$ws.Range("L135").Value = 750087000
$ws.Range("N135").Value = -750092070

$ws.Range("H136").Value = 60000
$ws.Range("J136").Value = 60000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -70200

$ws.Range("H138").Value = 2290.4062
$ws.Range("I138").Value = 623.93335
$ws.Range("J138").Value = 3760.8235
$ws.Range("K138").Value = 1871.80005
$ws.Range("L138").Value = 11282.4705
$ws.Range("M138").Value = 3268.19995
$ws.Range("N138").Value = -21562.4705


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1180
$ws.Range("I20").Value = 1033.3334
$ws.Range("K20").Value = 1033.3334
$ws.Range("M20").Value = -786.3334

$ws.Range("H99").Value = 802.95
$ws.Range("I99").Value = 791.7059
$ws.Range("J99").Value = 866.6667
$ws.Range("K99").Value = 791.7059
$ws.Range("L99").Value = 866.6667
$ws.Range("M99").Value = 706.2941
$ws.Range("N99").Value = -3862.6667


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9486.275
$ws.Range("J31").Value = 5738
$ws.Range("L31").Value = 5738
$ws.Range("N31").Value = -6328

$ws.Range("H34").Value = 9486.275
$ws.Range("J34").Value = 5738
$ws.Range("L34").Value = 5738
$ws.Range("N34").Value = -6142

$ws.Range("H58").Value = 27593.63
$ws.Range("I58").Value = 1322.2941
$ws.Range("J58").Value = 250900
$ws.Range("K58").Value = 1322.2941
$ws.Range("L58").Value = 250900
$ws.Range("M58").Value = -1119.2941
$ws.Range("N58").Value = -251306

$ws.Range("H62").Value = 4201.5
$ws.Range("I62").Value = 3457.1428
$ws.Range("J62").Value = 5243.6
$ws.Range("K62").Value = 3457.1428
$ws.Range("L62").Value = 5243.6
$ws.Range("M62").Value = -2833.1428
$ws.Range("N62").Value = -6491.6

$ws.Range("H65").Value = 4201.5
$ws.Range("I65").Value = 3457.1428
$ws.Range("J65").Value = 5243.6
$ws.Range("K65").Value = 17285.714
$ws.Range("L65").Value = 26218
$ws.Range("M65").Value = -14165.714
$ws.Range("N65").Value = -32458

$ws.Range("H86").Value = 12271.474
$ws.Range("J86").Value = 18628.666
$ws.Range("L86").Value = 18628.666
$ws.Range("N86").Value = -20874.666

$ws.Range("H89").Value = 12271.474
$ws.Range("J89").Value = 18628.666
$ws.Range("L89").Value = 93143.33
$ws.Range("N89").Value = -104375.33

$ws.Range("H136").Value = 27593.63
$ws.Range("I136").Value = 1322.2941
$ws.Range("J136").Value = 250900
$ws.Range("K136").Value = 3966.8823
$ws.Range("L136").Value = 752700
$ws.Range("M136").Value = -1416.8823
$ws.Range("N136").Value = -757800


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 750.1724
$ws.Range("I5").Value = 649.1111
$ws.Range("J5").Value = 795.65
$ws.Range("K5").Value = 1947.3333
$ws.Range("L5").Value = 2386.95
$ws.Range("M5").Value = -1835.3333
$ws.Range("N5").Value = -2610.95

$ws.Range("H37").Value = 75000
$ws.Range("J37").Value = 75000
$ws.Range("L37").Value = 225000
$ws.Range("N37").Value = -225224

$ws.Range("H129").Value = 500779.9
$ws.Range("I129").Value = 2000
$ws.Range("J129").Value = 556199.9
$ws.Range("K129").Value = 6000
$ws.Range("L129").Value = 1668599.7
$ws.Range("M129").Value = -1000
$ws.Range("N129").Value = -1678599.7

$ws.Range("H131").Value = 796.1900000000001
$ws.Range("J131").Value = 826.129
$ws.Range("L131").Value = 2478.387
$ws.Range("N131").Value = -12558.387

$ws.Range("H135").Value = 750.1724
$ws.Range("I135").Value = 649.1111
$ws.Range("J135").Value = 795.65
$ws.Range("K135").Value = 5841.9999
$ws.Range("L135").Value = 7160.849999999999
$ws.Range("M135").Value = -3306.9999
$ws.Range("N135").Value = -12230.85

$ws.Range("H140").Value = 3956.558
$ws.Range("I140").Value = 4632.269
$ws.Range("K140").Value = 13896.807
$ws.Range("M140").Value = -8716.807000000001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3429.0476
$ws.Range("I80").Value = 2800.7693
$ws.Range("J80").Value = 4450
$ws.Range("K80").Value = 2800.7693
$ws.Range("L80").Value = 4450
$ws.Range("M80").Value = -1802.7693
$ws.Range("N80").Value = -6446

$ws.Range("H83").Value = 3429.0476
$ws.Range("I83").Value = 2800.7693
$ws.Range("J83").Value = 4450
$ws.Range("K83").Value = 14003.8465
$ws.Range("L83").Value = 22250
$ws.Range("M83").Value = -9011.8465
$ws.Range("N83").Value = -32234

$ws.Range("H97").Value = 1598.4286
$ws.Range("I97").Value = 867.25
$ws.Range("K97").Value = 867.25
$ws.Range("M97").Value = -371.25


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1073.2222
$ws.Range("I22").Value = 563.3333
$ws.Range("J22").Value = 1328.1666
$ws.Range("K22").Value = 563.3333
$ws.Range("L22").Value = 1328.1666
$ws.Range("M22").Value = -268.3333
$ws.Range("N22").Value = -1918.1666

$ws.Range("H27").Value = 1073.2222
$ws.Range("I27").Value = 563.3333
$ws.Range("J27").Value = 1328.1666
$ws.Range("K27").Value = 563.3333
$ws.Range("L27").Value = 1328.1666
$ws.Range("M27").Value = -456.3333
$ws.Range("N27").Value = -1542.1666

$ws.Range("H55").Value = 259
$ws.Range("I55").Value = 164
$ws.Range("K55").Value = 164
$ws.Range("M55").Value = 9

$ws.Range("H110").Value = 36349.75
$ws.Range("J110").Value = 36349.75
$ws.Range("L110").Value = 36349.75
$ws.Range("N110").Value = -44529.75

